$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H4: new topic cell in the row-4 group, same formatting as the other G4-style cells
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = "Design Principles"

# H7: new topic cell in the row-7 group, same formatting as G7
$ws.Range("G7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = "Datastructures & Algorithms"

# H10: new topic cell in the row-10 group, same formatting as G10
$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Value = "Microservices "

# H13: new topic cell in the row-13 group, same formatting as G13
$ws.Range("G13").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = "Spring Cloud"

# Update the sheet view to reflect the new selection / scroll position
$ws.Range("H4").Select()
$excel.ActiveWindow.ScrollColumn = 4
